$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows before the existing row 870, shifting the
# rest of the table (old rows 870-907) down to rows 873-910.
$ws.Rows("870:872").Insert()

# Shared template values for this data block (same on every row of this
# sheet's data region).
$mercadoId = 10
$mercado   = "Vega Modelo de Temuco"
$region    = "La Araucanía"
$codreg    = 9
$tipo      = "Fruta"
$prodId    = 100108
$producto  = "Tropicales y subtropicales"
$catId     = 100108006
$categoria = "Plátano"
$unidad    = "`$/caja 20 kilos"
$origen    = "Ecuador"
$kgUnidad  = 20

function Set-Row {
    param($row, $fecha, $variedad, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $precioKg)

    $ws.Cells($row, 1).Value  = $mercadoId
    $ws.Cells($row, 2).Value  = $mercado
    $ws.Cells($row, 3).Value  = $region
    $ws.Cells($row, 4).Value  = $fecha
    $ws.Cells($row, 5).Value  = $codreg
    $ws.Cells($row, 6).Value  = $tipo
    $ws.Cells($row, 7).Value  = $prodId
    $ws.Cells($row, 8).Value  = $producto
    $ws.Cells($row, 9).Value  = $catId
    $ws.Cells($row, 10).Value = $categoria
    $ws.Cells($row, 11).Value = $variedad
    $ws.Cells($row, 12).Value = $calidad
    $ws.Cells($row, 13).Value = $volumen
    $ws.Cells($row, 14).Value = $precioMin
    $ws.Cells($row, 15).Value = $precioMax
    $ws.Cells($row, 16).Value = $precioProm
    $ws.Cells($row, 17).Value = $unidad
    $ws.Cells($row, 18).Value = $origen
    $ws.Cells($row, 19).Value = $precioKg
    $ws.Cells($row, 20).Value = $kgUnidad
}

Set-Row 870 44939 "Barraganete"     "Maduro" 40  40000 40000 40000 2000
Set-Row 871 44939 "Barraganete"     "Verde"  50  36000 36000 36000 1800
Set-Row 872 44939 "Sin especificar" "Pintón" 700 26000 27000 26571 1329
